# Natmi following Dr Hou advice
# The LR-pair table (Efna5 -> Epha7) is expanded from a 2x2 sending/target
# cluster grid (FAPs, sCs) to a full 3x3 grid that also includes the ECs
# cluster, recomputing all downstream NATMI specificity/weight metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ligand/receptor are constant across all rows of this sheet.
$ligand = "Efna5"
$receptor = "Epha7"

# Per-cluster values (Ligand-expressing cells, Ligand detection rate,
# Ligand average/total expression, Ligand derived specificity avg/total)
# keyed by cluster name - these only depend on the sending cluster.
$ligandStats = @{
    "ECs"  = @(2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572)
    "FAPs" = @(3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045)
    "sCs"  = @(3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798)
}

# Per-cluster values (Receptor-expressing cells, Receptor detection rate,
# Receptor average/total expression, Receptor derived specificity avg/total)
# keyed by cluster name - these only depend on the target cluster.
$receptorStats = @{
    "ECs"  = @(3, 1, 0.3173666666666666, 0.9520999999999999, 0.1863268828340317, 0.1863268828340316)
    "FAPs" = @(2, 0.6666666666666666, 0.2780386666666667, 0.8341160000000001, 0.1632373009158609, 0.1632373009158609)
    "sCs"  = @(3, 1, 1.107873666666667, 3.323621, 0.6504358162501074, 0.6504358162501074)
}

# Edge metrics (average weight, total weight, average specificity, total
# specificity) depend on the sending/target cluster pair.
$edgeStats = @{
    "ECs|ECs"   = @(0.07326568183333333, 0.6593911364999999,   0.01188263375659482, 0.01188263375659481)
    "ECs|FAPs"  = @(0.06418661639333334, 0.5776795475400001,   0.01041014067694133, 0.01041014067694133)
    "ECs|sCs"   = @(0.2557581753183333,  2.301823577865,       0.04148027632467957, 0.04148027632467957)
    "FAPs|ECs"  = @(0.9205954341555554,  8.285358907399999,    0.1493072623953519,  0.1493072623953519)
    "FAPs|FAPs" = @(0.8065154722782224,  7.258639250504001,    0.1308051428213017,  0.1308051428213017)
    "FAPs|sCs"  = @(3.213643858274889,   28.922794724474,      0.5212065463183508,  0.5212065463183509)
    "sCs|ECs"   = @(0.1549890795444444,  1.3949017159,         0.02513698668208498, 0.02513698668208497)
    "sCs|FAPs"  = @(0.1357828705737778,  1.222045835164,       0.02202201741761789, 0.02202201741761789)
    "sCs|sCs"   = @(0.5410408145621112,  4.869367331059,       0.08774899360707694, 0.08774899360707694)
}

$clusters = @("ECs", "FAPs", "sCs")

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $lstat = $ligandStats[$sending]
        $rstat = $receptorStats[$target]
        $estat = $edgeStats["$sending|$target"]

        $ws.Cells.Item($row, 1).Value  = $sending
        $ws.Cells.Item($row, 2).Value  = $ligand
        $ws.Cells.Item($row, 3).Value  = $receptor
        $ws.Cells.Item($row, 4).Value  = $target

        $ws.Cells.Item($row, 5).Value  = $lstat[0]
        $ws.Cells.Item($row, 6).Value  = $lstat[1]
        $ws.Cells.Item($row, 7).Value  = $lstat[2]
        $ws.Cells.Item($row, 8).Value  = $lstat[3]
        $ws.Cells.Item($row, 9).Value  = $lstat[4]
        $ws.Cells.Item($row, 10).Value = $lstat[5]

        $ws.Cells.Item($row, 11).Value = $rstat[0]
        $ws.Cells.Item($row, 12).Value = $rstat[1]
        $ws.Cells.Item($row, 13).Value = $rstat[2]
        $ws.Cells.Item($row, 14).Value = $rstat[3]
        $ws.Cells.Item($row, 15).Value = $rstat[4]
        $ws.Cells.Item($row, 16).Value = $rstat[5]

        $ws.Cells.Item($row, 17).Value = $estat[0]
        $ws.Cells.Item($row, 18).Value = $estat[1]
        $ws.Cells.Item($row, 19).Value = $estat[2]
        $ws.Cells.Item($row, 20).Value = $estat[3]

        $row++
    }
}
